$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number formats/styles from column E into the new column D
# (done per contiguous data block so we do not create stray cells
# in rows that have no data, e.g. 5, 6, 37, 79)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 678900
$ws.Range("D9").Value = 111900
$ws.Range("D10").Value = 567000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 202400
$ws.Range("D15").Value = 119800
$ws.Range("D17").Value = 448700
$ws.Range("D18").Value = 230200
$ws.Range("D20").Value = -600
$ws.Range("D21").Value = 349400
$ws.Range("D22").Value = 86000
$ws.Range("D23").Value = 143600
$ws.Range("D24").Value = -100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 143700
$ws.Range("D27").Value = 143700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 600
$ws.Range("D33").Value = 143700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 143700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 2400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 97600
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 128500
$ws.Range("D46").Value = 228400
$ws.Range("D47").Value = 5100
$ws.Range("D48").Value = 1202700
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 67400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1503600
$ws.Range("D57").Value = 135500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 96000
$ws.Range("D60").Value = 231500
$ws.Range("D61").Value = 830200
$ws.Range("D62").Value = 12100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1073800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -796900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 429900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 143700
$ws.Range("D83").Value = 119800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 244300
$ws.Range("D91").Value = -216800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -474500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 130400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -99800

# A handful of restated prior-period values (now shifted into E/F)
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = -1100
$ws.Range("E21").Value = 119000
$ws.Range("F21").Value = -169200
$ws.Range("E22").Value = 70300
$ws.Range("F22").Value = 64500
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1100

